$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Jakub Czulak column (H/I/J) updates for rows 11-13 ---
$ws.Range("J11").Value = 3
$ws.Range("J12").Value = 36
$ws.Range("J13").Value = 10

# --- New entries added to rows 14-16 (H/I/J columns) ---
# Copy date formatting from an existing date cell so the new cells pick up
# the same number format / style (s="3") instead of minting a new style.
$ws.Range("H13").Copy($ws.Range("H14"))
$ws.Range("H14").Value = 45770
$ws.Range("I14").Value = "PacjentServiceImpl.cs"
$ws.Range("J14").Value = 20

$ws.Range("H13").Copy($ws.Range("H15"))
$ws.Range("H15").Value = 45770
$ws.Range("I15").Value = "IPacjentService.cs"
$ws.Range("J15").Value = 3

$ws.Range("H13").Copy($ws.Range("H16"))
$ws.Range("H16").Value = 45770
$ws.Range("I16").Value = "Program.cs"
$ws.Range("J16").Value = 1

# --- View state: window scrolled so column F is the leftmost visible
# column, and the selection moves to K7 ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K7").Select()
